$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 21; $r -le 26; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "e.init") {
        $cell.Value = "e.close"
    }
}
